$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new data rows 9-31
$ws.Range("A9").Value = "HUSSAIN SAYYED"
$ws.Range("B9").Value = "hussain@gmail.com"
$ws.Range("C9").Value = 3803554
$ws.Range("D9").Value = "TYBBACA"
$ws.Range("E9").Value = 6618
$ws.Hyperlinks.Add($ws.Range("B9"), "mailto:hussain@gmail.com")

$ws.Range("A10").Value = "MUSAB DESHMUKH"
$ws.Range("B10").Value = "musab@gmail.com"
$ws.Range("C10").Value = 3801801
$ws.Range("D10").Value = "TYBBACA"
$ws.Range("E10").Value = 6605
$ws.Hyperlinks.Add($ws.Range("B10"), "mailto:musab@gmail.com")

$ws.Range("A11").Value = "SAQIB BASHIR DUSTE"
$ws.Range("B11").Value = "saquib@gmail.com"
$ws.Range("C11").Value = 3803349
$ws.Range("D11").Value = "TYBBACA"
$ws.Range("E11").Value = 6602
$ws.Hyperlinks.Add($ws.Range("B11"), "mailto:saquib@gmail.com")

$ws.Range("A12").Value = "AMAAN MAINODDIN SHAIKH"
$ws.Range("B12").Value = "amman@gmail.com"
$ws.Range("C12").Value = 3800982
$ws.Range("D12").Value = "TYBBACA"
$ws.Range("E12").Value = 6607
$ws.Hyperlinks.Add($ws.Range("B12"), "mailto:amman@gmail.com")

$ws.Range("A13").Value = "HAARIS AHMED MATEEN SHAIKH"
$ws.Range("B13").Value = "haaris@gmail.com"
$ws.Range("C13").Value = 3848678
$ws.Range("D13").Value = "TYBBACA"
$ws.Range("E13").Value = 6609
$ws.Hyperlinks.Add($ws.Range("B13"), "mailto:haaris@gmail.com")

$ws.Range("A14").Value = "ZIDAAN AASIF TAMBOLI"
$ws.Range("B14").Value = "zidan@gmail.com"
$ws.Range("C14").Value = 3995727
$ws.Range("D14").Value = "TYBBACA"
$ws.Range("E14").Value = 6610
$ws.Hyperlinks.Add($ws.Range("B14"), "mailto:zidan@gmail.com")

$ws.Range("A15").Value = "ASAD MOULALI SHAIKH"
$ws.Range("B15").Value = "naruto@gmail.com"
$ws.Range("C15").Value = 4813594
$ws.Range("D15").Value = "TYBBACA"
$ws.Range("E15").Value = 6613
$ws.Hyperlinks.Add($ws.Range("B15"), "mailto:naruto@gmail.com")

$ws.Range("A16").Value = "ADITYA PRAKASH SONULE"
$ws.Range("B16").Value = "ncc@gmail.com"
$ws.Range("C16").Value = 3800943
$ws.Range("D16").Value = "TYBBACA"
$ws.Range("E16").Value = 6625
$ws.Hyperlinks.Add($ws.Range("B16"), "mailto:ncc@gmail.com")

$ws.Range("A17").Value = "MUSKAN IRFAN PATWEKAR"
$ws.Range("B17").Value = "ceo@gmail.com"
$ws.Range("C17").Value = 3802610
$ws.Range("D17").Value = "TYBBACA"
$ws.Range("E17").Value = 6624
$ws.Hyperlinks.Add($ws.Range("B17"), "mailto:ceo@gmail.com")

$ws.Range("A18").Value = "ARYAN SACHIN GAIKWAD"
$ws.Range("B18").Value = "ayan@gmail.com"
$ws.Range("C18").Value = 3801327
$ws.Range("D18").Value = "TYBBACA"
$ws.Range("E18").Value = 6651
$ws.Hyperlinks.Add($ws.Range("B18"), "mailto:ayan@gmail.com")

$ws.Range("A19").Value = "AKIL SHAKIL FARAS"
$ws.Range("B19").Value = "akil@gamil.com"
$ws.Range("C19").Value = 3997262
$ws.Range("D19").Value = "TYBBACA"
$ws.Range("E19").Value = 6655
$ws.Hyperlinks.Add($ws.Range("B19"), "mailto:akil@gamil.com")

$ws.Range("A20").Value = "MD JUNAID ASHRAF SHAIKH"
$ws.Range("B20").Value = "junu@gamil.com"
$ws.Range("C20").Value = 3801331
$ws.Range("D20").Value = "TYBBACA"
$ws.Range("E20").Value = 6667
$ws.Hyperlinks.Add($ws.Range("B20"), "mailto:junu@gamil.com")

$ws.Range("A21").Value = "UDAY RAJESH SHINDE"
$ws.Range("B21").Value = "uday@gamil.com"
$ws.Range("C21").Value = 4261020
$ws.Range("D21").Value = "TYBBACA"
$ws.Range("E21").Value = 6679
$ws.Hyperlinks.Add($ws.Range("B21"), "mailto:uday@gamil.com")

$ws.Range("A22").Value = "HUSSAIN ABBAS SAYYED"
$ws.Range("B22").Value = "husssain@gamil.com"
$ws.Range("C22").Value = 3998954
$ws.Range("D22").Value = "TYBBACA"
$ws.Range("E22").Value = 6677
$ws.Hyperlinks.Add($ws.Range("B22"), "mailto:husssain@gamil.com")

$ws.Range("A23").Value = "RUHANA SARFRAZ JAHAGIRDAR"
$ws.Range("B23").Value = "ruhana@gamil.com"
$ws.Range("C23").Value = 3848309
$ws.Range("D23").Value = "TYBBACA"
$ws.Range("E23").Value = 6676
$ws.Hyperlinks.Add($ws.Range("B23"), "mailto:ruhana@gamil.com")

$ws.Range("A24").Value = "DEVENDRA HANUMANTH GAIKWAD"
$ws.Range("B24").Value = "dev@gamil.com"
$ws.Range("C24").Value = 4002311
$ws.Range("D24").Value = "TYBBACA"
$ws.Range("E24").Value = 6663
$ws.Hyperlinks.Add($ws.Range("B24"), "mailto:dev@gamil.com")

$ws.Range("A25").Value = "SAHIL ALI MUSHTAQUE ALI SHAIKH"
$ws.Range("B25").Value = "sahil@gamil.com"
$ws.Range("C25").Value = 3803431
$ws.Range("D25").Value = "TYBBACA"
$ws.Range("E25").Value = 6661
$ws.Hyperlinks.Add($ws.Range("B25"), "mailto:sahil@gamil.com")

$ws.Range("A26").Value = "SHOAIB SAMEER SHAIKH"
$ws.Range("B26").Value = "sammer@gamil.com"
$ws.Range("C26").Value = 3801362
$ws.Range("D26").Value = "TYBBACA"
$ws.Range("E26").Value = 6653
$ws.Hyperlinks.Add($ws.Range("B26"), "mailto:sammer@gamil.com")

$ws.Range("A27").Value = "ALTAF NABI KAKANDKI"
$ws.Range("B27").Value = "altaf@gamil.com"
$ws.Range("C27").Value = 3848745
$ws.Range("D27").Value = "TYBBACA"
$ws.Range("E27").Value = 6648
$ws.Hyperlinks.Add($ws.Range("B27"), "mailto:altaf@gamil.com")

$ws.Range("A28").Value = "MOHAMMED NATIQ HASAN SAYYED"
$ws.Range("B28").Value = "natiq@gamil.com"
$ws.Range("C28").Value = 3803350
$ws.Range("D28").Value = "TYBBACA"
$ws.Range("E28").Value = 6638
$ws.Hyperlinks.Add($ws.Range("B28"), "mailto:natiq@gamil.com")

$ws.Range("A29").Value = "TAAHA KHALID SIDDIQUI"
$ws.Range("B29").Value = "taha@gamil.com"
$ws.Range("C29").Value = 4012808
$ws.Range("D29").Value = "TYBBACA"
$ws.Range("E29").Value = 6635
$ws.Hyperlinks.Add($ws.Range("B29"), "mailto:taha@gamil.com")

$ws.Range("A30").Value = "RASAAM NAEEM BHALDAR"
$ws.Range("B30").Value = "rassam@gamil.com"
$ws.Range("C30").Value = 3800825
$ws.Range("D30").Value = "TYBBACA"
$ws.Range("E30").Value = 6634
$ws.Hyperlinks.Add($ws.Range("B30"), "mailto:rassam@gamil.com")

$ws.Range("A31").Value = "MUFIZ QUTBUDDIN SHAIKH"
$ws.Range("B31").Value = "mufiz@gamil.com"
$ws.Range("C31").Value = 4851066
$ws.Range("D31").Value = "TYBBACA"
$ws.Range("E31").Value = 6628
$ws.Hyperlinks.Add($ws.Range("B31"), "mailto:mufiz@gamil.com")

# Re-apply Hyperlink style across B9:B31 in one pass (keeps style table reuse minimal)
$ws.Range("B9:B31").Style = "Hyperlink"

# Column width changes
$ws.Columns("A").ColumnWidth = 28.33
$ws.Columns("B").ColumnWidth = 23.5

# View state: active cell B31, scrolled so row 9 is at top
[void]$ws.Range("B31").Select()
$excel.ActiveWindow.ScrollRow = 9
$excel.ActiveWindow.ScrollColumn = 1
